$p = $ppt.ActivePresentation

# Reposition/resize the comparison picture (Picture 11) on slide 2.
# PowerPoint COM exposes Left/Top/Width/Height in points (Single/float32),
# while the target values come from the OOXML in EMUs (1 pt = 12700 EMU).
# Add a half-EMU nudge before converting so the float32 round-trip lands
# on the exact target EMU value instead of truncating one unit short.
$s2 = $p.Slides.Item(2)
$pic = $s2.Shapes.Item(3)
$pic.Left = (6342380 + 0.5) / 12700
$pic.Top = (2721909 + 0.5) / 12700
$pic.Width = (5303533 + 0.5) / 12700
$pic.Height = (3030590 + 0.5) / 12700

# Remove the last slide (sldId 260 / slide4.xml, the "Final curve" slide).
$s4 = $p.Slides.Item(4)
$s4.Delete()
